# Generate Report for Handoff
#
# Refreshes the handoff timestamps produced by the report generator and
# marks the just-handed-off rows with their "ht" (handoff type) priority.
#
#  - Overview!G and de-de!H share the "2016-08-22 12:19:53" timestamp ->
#    "2016-08-22 12:20:31".
#  - zh-cn!H carries the "2016-08-22 12:19:48" timestamp ->
#    "2016-08-22 12:20:25".
#  - zh-cn!E and de-de!E (Priority) go from blank to "ht" for the rows that
#    were just handed off.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-22 12:20:31"
}

foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-22 12:20:25"
}

foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-22 12:20:31"
}
